$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select D2:E2 (the "Successful"/"Passed" assertion-result cells) and clear them out,
# matching the author deleting the test-run status values from the login data sheet.
$ws.Range("D2:E2").Select()
$ws.Range("D2:E2").ClearContents()
